# Consumer signin, consumer bookings
# Populate Sheet1 with consumer credential / access-token data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (character-width units; engine rounds to nearest 1/6) ---
$ws.Columns.Item(1).ColumnWidth = 25.333333333333336
$ws.Columns.Item(2).ColumnWidth = 14.166666666666666

# --- Header row ---
$ws.Range("A1").Value = "Email"
$ws.Range("B1").Value = "Password"

# --- Data row (email / password) ---
$ws.Range("A2").Value = "consumer2@gmail.com"
$ws.Range("B2").Value = "password@123"

# --- Access-token header + header for C1 must come after A1/B1 but the
#     value itself is independent of the other columns. ---
$ws.Range("C1").Value = "Access token"
$ws.Range("C2").Value = "eyJ0eXAiOiJKV1QiLCJhbGciOiJIUzI1NiJ9.eyJ0b2tlbl90eXBlIjoiYWNjZXNzIiwiZXhwIjoxNjUxMDc2NzgzLCJqdGkiOiJmNjA2YjliNGJlYWU0YThkYTdlNTczNzRmYjE3YTdmZiIsInVzZXJfaWQiOjY4NSwicm9sZSI6IkMiLCJ1c2VybmFtZSI6ImNvbnN1bWVyMiIsImVtYWlsIjoiY29uc3VtZXIyQGdtYWlsLmNvbSJ9.1lla5-4SwRGRNWRrZaBlii-135vVOJ4Z8AJ-IairJ8I"

# --- Hyperlinks on the email & password cells, pointing at the consumer's
#     access-token values captured from the signin responses. ---
[void]$ws.Hyperlinks.Add($ws.Range("A2"), "eyJ0eXAiOiJKV1QiLCJhbGciOiJIUzI1NiJ9.eyJ0b2tlbl90eXBlIjoiYWNjZXNzIiwiZXhwIjoxNjUxMDc1NjExLCJqdGkiOiI3YTQ1MjQ2YTQ5MGM0ZDMwOTljMzYzNjdhOGYxZmI1YiIsInVzZXJfaWQiOjY4NSwicm9sZSI6IkMiLCJ1c2VybmFtZSI6ImNvbnN1bWVyMiIsImVtYWlsIjoiY29uc3VtZXIyQGdtYWlsLmNvbSJ9.4QB2bD7_-qmT5Msq6qIqynrK6iKBVAugjRlhFUpR2Hg")
[void]$ws.Hyperlinks.Add($ws.Range("B2"), "eyJ0eXAiOiJKV1QiLCJhbGciOiJIUzI1NiJ9.eyJ0b2tlbl90eXBlIjoiYWNjZXNzIiwiZXhwIjoxNjUxMDc1NzcyLCJqdGkiOiI4ZTJkZmQ2M2MzNzQ0NmJiYWMzNDZlMzcyYzNhZWFjZCIsInVzZXJfaWQiOjY4NSwicm9sZSI6IkMiLCJ1c2VybmFtZSI6ImNvbnN1bWVyMiIsImVtYWlsIjoiY29uc3VtZXIyQGdtYWlsLmNvbSJ9.Z3DLms6ucRIdiN8Oad285ZyBFUUO4d0u53pSfa4KVQM")

# --- Page setup: portrait, A4 ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Leave the active selection on C1, matching the authored workbook ---
[void]$ws.Range("C1").Select()
